$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.575.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.666.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9977'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4630'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.69%  '
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06137'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.662.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.336'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '75.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5727'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9992'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.569.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006713'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.873.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.421'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.11%  '
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.219'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.372'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.718'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.945'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07680'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.594'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04343'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6059'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9396'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9169'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '107.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.357'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9981'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.834'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.15%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.023'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.58%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3713'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05257'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.589'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9996'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9982'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.02%  '
